$d = $word.ActiveDocument

$d.Content.Find.Execute("IV catheter placed into a vein in the hand or arm", $true, $false, $false, $false, $false, $true, 1, $false, "IV catheter placed in vein of hand or arm", 2)
$d.Content.Find.Execute("Allows administration of chemotherapy and fluids", $true, $false, $false, $false, $false, $true, 1, $false, "Allows administration of chemo and fluids", 2)
$d.Content.Find.Execute("Not suitable for FLOT chemotherapy", $true, $false, $false, $false, $false, $true, 1, $false, "Not suitable for FLOT chemo", 2)
$d.Content.Find.Execute("Implantable device that makes the administration of chemotherapy easier", $true, $false, $false, $false, $false, $true, 1, $false, "Implantable device makes chemo easier", 2)
$d.Content.Find.Execute("Suitable for FLOT chemo", $true, $true, $false, $false, $false, $true, 1, $false, "OK for FLOT chemo", 2)
$d.Content.Find.Execute("Sutures dissolve on their own", $true, $false, $false, $false, $false, $true, 1, $false, "Sutures dissolve", 2)
